$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "x" progress cells for the first weeks already completed
$ws.Range("D4").Value = "x"
$ws.Range("D5").Value = "x"
$ws.Range("E6").Value = "x"
$ws.Range("F7").Value = "x"
$ws.Range("F8").Value = "x"
$ws.Range("G9").Value = "x"

# Update activity descriptions (add final period / fix accent, add new row)
$ws.Range("B7").Value = "Identificar metodologia existente para categorizar el nivel de seguridad en la organización."
$ws.Range("B8").Value = "Identificar metodologia existente para seleccionar los controles de seguridad."
$ws.Range("B9").Value = "Identificar procedimientos y herramientas para dar el seguimiento de las implementaciones de los controles de seguridad."
$ws.Range("B10").Value = "Establecer formatos de encuestas para que los arquitectos de software puedan completar."

# Mark the new activity as in-progress and underline the mark
$ws.Range("G10").Value = "x"
$ws.Range("G10").Font.Underline = $true

$ws.Range("G10").Select()
